$wb = $excel.ActiveWorkbook

# ===== Sheet 1 =====
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(1, 1).Value = 'LÍNEA 141 - LP1912 - 24/01/2026'
$ws.Cells.Item(2, 1).Value = 'Última actualización: 10:51:29'
$ws.Cells.Item(3, 1).Value = 'Total filas: 182'

$arr = New-Object "object[,]" 183,5
$arr[0,0] = 'Hora_Scrap'
$arr[0,1] = 'Hora_Llegada'
$arr[0,2] = 'Linea'
$arr[0,3] = 'Minutos'
$arr[0,4] = 'Parada'
$arr[1,0] = '03:52:04'
$arr[1,1] = '04:01'
$arr[1,2] = '81_EL PELIGRO'
$arr[1,3] = 9
$arr[1,4] = 'LP1912'
$arr[2,0] = '04:32:18'
$arr[2,1] = '04:33'
$arr[2,2] = '15_ABASTO'
$arr[2,3] = 1
$arr[2,4] = 'LP1912'
$arr[3,0] = '03:52:04'
$arr[3,1] = '04:46'
$arr[3,2] = '215A_EL PATO'
$arr[3,3] = 54
$arr[3,4] = 'LP1912'
$arr[4,0] = '04:44:46'
$arr[4,1] = '04:46'
$arr[4,2] = '215_EL PELIGRO'
$arr[4,3] = 2
$arr[4,4] = 'LP1912'
$arr[5,0] = '04:44:46'
$arr[5,1] = '04:46'
$arr[5,2] = '15_ABASTO'
$arr[5,3] = 2
$arr[5,4] = 'LP1912'
$arr[6,0] = '04:32:18'
$arr[6,1] = '04:47'
$arr[6,2] = '215_EL PELIGRO'
$arr[6,3] = 15
$arr[6,4] = 'LP1912'
$arr[7,0] = '04:52:24'
$arr[7,1] = '04:53'
$arr[7,2] = '11_ETCHEVERRY'
$arr[7,3] = 1
$arr[7,4] = 'LP1912'
$arr[8,0] = '04:52:24'
$arr[8,1] = '04:54'
$arr[8,2] = '15_ABASTO'
$arr[8,3] = 2
$arr[8,4] = 'LP1912'
$arr[9,0] = '04:13:31'
$arr[9,1] = '05:11'
$arr[9,2] = '17_ROMERO'
$arr[9,3] = 58
$arr[9,4] = 'LP1912'
$arr[10,0] = '03:52:04'
$arr[10,1] = '05:16'
$arr[10,2] = '17_ROMERO'
$arr[10,3] = 84
$arr[10,4] = 'LP1912'
$arr[11,0] = '04:52:24'
$arr[11,1] = '05:22'
$arr[11,2] = '23_HERNANDEZ'
$arr[11,3] = 30
$arr[11,4] = 'LP1912'
$arr[12,0] = '05:16:02'
$arr[12,1] = '05:25'
$arr[12,2] = '23_HERNANDEZ'
$arr[12,3] = 9
$arr[12,4] = 'LP1912'
$arr[13,0] = '04:44:46'
$arr[13,1] = '05:31'
$arr[13,2] = '81_EL PELIGRO'
$arr[13,3] = 47
$arr[13,4] = 'LP1912'
$arr[14,0] = '05:16:02'
$arr[14,1] = '05:32'
$arr[14,2] = '81_EL PELIGRO'
$arr[14,3] = 16
$arr[14,4] = 'LP1912'
$arr[15,0] = '03:52:04'
$arr[15,1] = '05:35'
$arr[15,2] = '215B_EL PATO'
$arr[15,3] = 103
$arr[15,4] = 'LP1912'
$arr[16,0] = '05:43:29'
$arr[16,1] = '05:44'
$arr[16,2] = '14_ABASTO'
$arr[16,3] = 1
$arr[16,4] = 'LP1912'
$arr[17,0] = '03:52:04'
$arr[17,1] = '05:46'
$arr[17,2] = '15_ABASTO'
$arr[17,3] = 114
$arr[17,4] = 'LP1912'
$arr[18,0] = '04:32:18'
$arr[18,1] = '05:47'
$arr[18,2] = '14_ABASTO'
$arr[18,3] = 75
$arr[18,4] = 'LP1912'
$arr[19,0] = '04:13:31'
$arr[19,1] = '05:50'
$arr[19,2] = '14_ABASTO'
$arr[19,3] = 97
$arr[19,4] = 'LP1912'
$arr[20,0] = '04:44:46'
$arr[20,1] = '05:51'
$arr[20,2] = '17_ROMERO'
$arr[20,3] = 67
$arr[20,4] = 'LP1912'
$arr[21,0] = '05:43:29'
$arr[21,1] = '05:52'
$arr[21,2] = '17_ROMERO'
$arr[21,3] = 9
$arr[21,4] = 'LP1912'
$arr[22,0] = '05:57:38'
$arr[22,1] = '06:00'
$arr[22,2] = '16_SANTA ANA'
$arr[22,3] = 3
$arr[22,4] = 'LP1912'
$arr[23,0] = '05:43:29'
$arr[23,1] = '06:01'
$arr[23,2] = '16_SANTA ANA'
$arr[23,3] = 18
$arr[23,4] = 'LP1912'
$arr[24,0] = '05:57:38'
$arr[24,1] = '06:03'
$arr[24,2] = '10_OLMOS'
$arr[24,3] = 6
$arr[24,4] = 'LP1912'
$arr[25,0] = '05:16:02'
$arr[25,1] = '06:04'
$arr[25,2] = '10_OLMOS'
$arr[25,3] = 48
$arr[25,4] = 'LP1912'
$arr[26,0] = '05:57:38'
$arr[26,1] = '06:10'
$arr[26,2] = '215A_EL PATO'
$arr[26,3] = 13
$arr[26,4] = 'LP1912'
$arr[27,0] = '05:43:29'
$arr[27,1] = '06:11'
$arr[27,2] = '215A_EL PATO'
$arr[27,3] = 28
$arr[27,4] = 'LP1912'
$arr[28,0] = '04:32:18'
$arr[28,1] = '06:15'
$arr[28,2] = '17_ROMERO'
$arr[28,3] = 103
$arr[28,4] = 'LP1912'
$arr[29,0] = '06:18:01'
$arr[29,1] = '06:24'
$arr[29,2] = '11_ETCHEVERRY'
$arr[29,3] = 6
$arr[29,4] = 'LP1912'
$arr[30,0] = '06:18:01'
$arr[30,1] = '06:27'
$arr[30,2] = '23_HERNANDEZ'
$arr[30,3] = 9
$arr[30,4] = 'LP1912'
$arr[31,0] = '04:44:46'
$arr[31,1] = '06:28'
$arr[31,2] = '17_ROMERO'
$arr[31,3] = 104
$arr[31,4] = 'LP1912'
$arr[32,0] = '04:44:46'
$arr[32,1] = '06:30'
$arr[32,2] = '16_SANTA ANA'
$arr[32,3] = 106
$arr[32,4] = 'LP1912'
$arr[33,0] = '06:18:01'
$arr[33,1] = '06:31'
$arr[33,2] = '17X38_ROMERO'
$arr[33,3] = 13
$arr[33,4] = 'LP1912'
$arr[34,0] = '06:18:01'
$arr[34,1] = '06:31'
$arr[34,2] = '16_SANTA ANA'
$arr[34,3] = 13
$arr[34,4] = 'LP1912'
$arr[35,0] = '04:52:24'
$arr[35,1] = '06:36'
$arr[35,2] = '17_ROMERO'
$arr[35,3] = 104
$arr[35,4] = 'LP1912'
$arr[36,0] = '06:35:33'
$arr[36,1] = '06:39'
$arr[36,2] = '225_C ROCA-H SUR'
$arr[36,3] = 4
$arr[36,4] = 'LP1912'
$arr[37,0] = '06:35:33'
$arr[37,1] = '06:41'
$arr[37,2] = '17_ROMERO'
$arr[37,3] = 6
$arr[37,4] = 'LP1912'
$arr[38,0] = '06:18:01'
$arr[38,1] = '06:45'
$arr[38,2] = '17_ROMERO'
$arr[38,3] = 27
$arr[38,4] = 'LP1912'
$arr[39,0] = '06:46:37'
$arr[39,1] = '06:46'
$arr[39,2] = '17_ROMERO'
$arr[39,3] = 0
$arr[39,4] = 'LP1912'
$arr[40,0] = '06:46:37'
$arr[40,1] = '06:50'
$arr[40,2] = '215A_EL PATO'
$arr[40,3] = 4
$arr[40,4] = 'LP1912'
$arr[41,0] = '05:16:02'
$arr[41,1] = '06:50'
$arr[41,2] = '17_ROMERO'
$arr[41,3] = 94
$arr[41,4] = 'LP1912'
$arr[42,0] = '06:35:33'
$arr[42,1] = '06:51'
$arr[42,2] = '215A_EL PATO'
$arr[42,3] = 16
$arr[42,4] = 'LP1912'
$arr[43,0] = '06:53:56'
$arr[43,1] = '06:53'
$arr[43,2] = '14_ABASTO'
$arr[43,3] = 0
$arr[43,4] = 'LP1912'
$arr[44,0] = '06:53:56'
$arr[44,1] = '06:54'
$arr[44,2] = '17_ROMERO'
$arr[44,3] = 1
$arr[44,4] = 'LP1912'
$arr[45,0] = '06:46:37'
$arr[45,1] = '06:54'
$arr[45,2] = '14_ABASTO'
$arr[45,3] = 8
$arr[45,4] = 'LP1912'
$arr[46,0] = '06:53:56'
$arr[46,1] = '07:03'
$arr[46,2] = '225_GOMEZ'
$arr[46,3] = 10
$arr[46,4] = 'LP1912'
$arr[47,0] = '06:46:37'
$arr[47,1] = '07:04'
$arr[47,2] = '225_GOMEZ'
$arr[47,3] = 18
$arr[47,4] = 'LP1912'
$arr[48,0] = '06:53:56'
$arr[48,1] = '07:06'
$arr[48,2] = '215C_EL PATO'
$arr[48,3] = 13
$arr[48,4] = 'LP1912'
$arr[49,0] = '06:18:01'
$arr[49,1] = '07:07'
$arr[49,2] = '215C_EL PATO'
$arr[49,3] = 49
$arr[49,4] = 'LP1912'
$arr[50,0] = '07:12:47'
$arr[50,1] = '07:13'
$arr[50,2] = '14X44_ABASTO'
$arr[50,3] = 1
$arr[50,4] = 'LP1912'
$arr[51,0] = '06:18:01'
$arr[51,1] = '07:14'
$arr[51,2] = '14X44_ABASTO'
$arr[51,3] = 56
$arr[51,4] = 'LP1912'
$arr[52,0] = '07:12:47'
$arr[52,1] = '07:20'
$arr[52,2] = '215A_EL PATO'
$arr[52,3] = 8
$arr[52,4] = 'LP1912'
$arr[53,0] = '06:35:33'
$arr[53,1] = '07:21'
$arr[53,2] = '215A_EL PATO'
$arr[53,3] = 46
$arr[53,4] = 'LP1912'
$arr[54,0] = '07:12:47'
$arr[54,1] = '07:23'
$arr[54,2] = '16_SANTA ANA'
$arr[54,3] = 11
$arr[54,4] = 'LP1912'
$arr[55,0] = '06:46:37'
$arr[55,1] = '07:24'
$arr[55,2] = '16_SANTA ANA'
$arr[55,3] = 38
$arr[55,4] = 'LP1912'
$arr[56,0] = '06:53:56'
$arr[56,1] = '07:28'
$arr[56,2] = '14_ABASTO'
$arr[56,3] = 35
$arr[56,4] = 'LP1912'
$arr[57,0] = '07:12:47'
$arr[57,1] = '07:29'
$arr[57,2] = '14_ABASTO'
$arr[57,3] = 17
$arr[57,4] = 'LP1912'
$arr[58,0] = '07:12:47'
$arr[58,1] = '07:33'
$arr[58,2] = '23_HERNANDEZ'
$arr[58,3] = 21
$arr[58,4] = 'LP1912'
$arr[59,0] = '07:12:47'
$arr[59,1] = '07:35'
$arr[59,2] = '17X38_ROMERO'
$arr[59,3] = 23
$arr[59,4] = 'LP1912'
$arr[60,0] = '06:46:37'
$arr[60,1] = '07:36'
$arr[60,2] = '17X38_ROMERO'
$arr[60,3] = 50
$arr[60,4] = 'LP1912'
$arr[61,0] = '07:12:47'
$arr[61,1] = '07:36'
$arr[61,2] = '27_EL RETIRO'
$arr[61,3] = 24
$arr[61,4] = 'LP1912'
$arr[62,0] = '06:18:01'
$arr[62,1] = '07:37'
$arr[62,2] = '27_EL RETIRO'
$arr[62,3] = 79
$arr[62,4] = 'LP1912'
$arr[63,0] = '07:12:47'
$arr[63,1] = '07:41'
$arr[63,2] = '16_SANTA ANA'
$arr[63,3] = 29
$arr[63,4] = 'LP1912'
$arr[64,0] = '07:12:47'
$arr[64,1] = '07:43'
$arr[64,2] = '10_OLMOS'
$arr[64,3] = 31
$arr[64,4] = 'LP1912'
$arr[65,0] = '06:18:01'
$arr[65,1] = '07:44'
$arr[65,2] = '10_OLMOS'
$arr[65,3] = 86
$arr[65,4] = 'LP1912'
$arr[66,0] = '07:38:30'
$arr[66,1] = '07:49'
$arr[66,2] = '15_ABASTO'
$arr[66,3] = 11
$arr[66,4] = 'LP1912'
$arr[67,0] = '07:50:33'
$arr[67,1] = '07:50'
$arr[67,2] = '15_ABASTO'
$arr[67,3] = 0
$arr[67,4] = 'LP1912'
$arr[68,0] = '07:50:33'
$arr[68,1] = '07:56'
$arr[68,2] = '10_OLMOS'
$arr[68,3] = 6
$arr[68,4] = 'LP1912'
$arr[69,0] = '07:50:33'
$arr[69,1] = '07:58'
$arr[69,2] = '23_HERNANDEZ'
$arr[69,3] = 8
$arr[69,4] = 'LP1912'
$arr[70,0] = '07:50:33'
$arr[70,1] = '07:59'
$arr[70,2] = '11_ETCHEVERRY'
$arr[70,3] = 9
$arr[70,4] = 'LP1912'
$arr[71,0] = '06:53:56'
$arr[71,1] = '07:59'
$arr[71,2] = '23_HERNANDEZ'
$arr[71,3] = 66
$arr[71,4] = 'LP1912'
$arr[72,0] = '06:18:01'
$arr[72,1] = '08:00'
$arr[72,2] = '11_ETCHEVERRY'
$arr[72,3] = 102
$arr[72,4] = 'LP1912'
$arr[73,0] = '06:46:37'
$arr[73,1] = '08:00'
$arr[73,2] = '23_HERNANDEZ'
$arr[73,3] = 74
$arr[73,4] = 'LP1912'
$arr[74,0] = '07:50:33'
$arr[74,1] = '08:01'
$arr[74,2] = '16_SANTA ANA'
$arr[74,3] = 11
$arr[74,4] = 'LP1912'
$arr[75,0] = '07:50:33'
$arr[75,1] = '08:03'
$arr[75,2] = '17X38_ROMERO'
$arr[75,3] = 13
$arr[75,4] = 'LP1912'
$arr[76,0] = '06:53:56'
$arr[76,1] = '08:13'
$arr[76,2] = '10_OLMOS'
$arr[76,3] = 80
$arr[76,4] = 'LP1912'
$arr[77,0] = '07:50:33'
$arr[77,1] = '08:14'
$arr[77,2] = '10_OLMOS'
$arr[77,3] = 24
$arr[77,4] = 'LP1912'
$arr[78,0] = '08:10:38'
$arr[78,1] = '08:19'
$arr[78,2] = '17_ROMERO'
$arr[78,3] = 9
$arr[78,4] = 'LP1912'
$arr[79,0] = '08:10:38'
$arr[79,1] = '08:21'
$arr[79,2] = '16_SANTA ANA'
$arr[79,3] = 11
$arr[79,4] = 'LP1912'
$arr[80,0] = '08:29:58'
$arr[80,1] = '08:29'
$arr[80,2] = '14_ABASTO'
$arr[80,3] = 0
$arr[80,4] = 'LP1912'
$arr[81,0] = '08:10:38'
$arr[81,1] = '08:33'
$arr[81,2] = '23_HERNANDEZ'
$arr[81,3] = 23
$arr[81,4] = 'LP1912'
$arr[82,0] = '08:29:58'
$arr[82,1] = '08:33'
$arr[82,2] = '215C_EL PATO'
$arr[82,3] = 4
$arr[82,4] = 'LP1912'
$arr[83,0] = '07:50:33'
$arr[83,1] = '08:34'
$arr[83,2] = '215C_EL PATO'
$arr[83,3] = 44
$arr[83,4] = 'LP1912'
$arr[84,0] = '08:29:58'
$arr[84,1] = '08:41'
$arr[84,2] = '16_SANTA ANA'
$arr[84,3] = 12
$arr[84,4] = 'LP1912'
$arr[85,0] = '08:29:58'
$arr[85,1] = '08:45'
$arr[85,2] = '10_OLMOS'
$arr[85,3] = 16
$arr[85,4] = 'LP1912'
$arr[86,0] = '08:40:53'
$arr[86,1] = '08:47'
$arr[86,2] = '215A_EL PATO'
$arr[86,3] = 7
$arr[86,4] = 'LP1912'
$arr[87,0] = '08:40:53'
$arr[87,1] = '08:47'
$arr[87,2] = '10_OLMOS'
$arr[87,3] = 7
$arr[87,4] = 'LP1912'
$arr[88,0] = '08:10:38'
$arr[88,1] = '08:48'
$arr[88,2] = '215A_EL PATO'
$arr[88,3] = 38
$arr[88,4] = 'LP1912'
$arr[89,0] = '08:10:38'
$arr[89,1] = '08:48'
$arr[89,2] = '10_OLMOS'
$arr[89,3] = 38
$arr[89,4] = 'LP1912'
$arr[90,0] = '08:29:58'
$arr[90,1] = '08:50'
$arr[90,2] = '16_P MOR-SANTA ANA'
$arr[90,3] = 21
$arr[90,4] = 'LP1912'
$arr[91,0] = '08:40:53'
$arr[91,1] = '08:51'
$arr[91,2] = '16_P MOR-SANTA ANA'
$arr[91,3] = 11
$arr[91,4] = 'LP1912'
$arr[92,0] = '08:52:13'
$arr[92,1] = '08:52'
$arr[92,2] = '16_P MOR-SANTA ANA'
$arr[92,3] = 0
$arr[92,4] = 'LP1912'
$arr[93,0] = '08:40:53'
$arr[93,1] = '08:59'
$arr[93,2] = '215B_EL PATO'
$arr[93,3] = 19
$arr[93,4] = 'LP1912'
$arr[94,0] = '08:10:38'
$arr[94,1] = '09:00'
$arr[94,2] = '23_HERNANDEZ'
$arr[94,3] = 50
$arr[94,4] = 'LP1912'
$arr[95,0] = '08:52:13'
$arr[95,1] = '09:00'
$arr[95,2] = '215B_EL PATO'
$arr[95,3] = 8
$arr[95,4] = 'LP1912'
$arr[96,0] = '08:52:13'
$arr[96,1] = '09:01'
$arr[96,2] = '16_SANTA ANA'
$arr[96,3] = 9
$arr[96,4] = 'LP1912'
$arr[97,0] = '07:38:30'
$arr[97,1] = '09:02'
$arr[97,2] = '23_HERNANDEZ'
$arr[97,3] = 84
$arr[97,4] = 'LP1912'
$arr[98,0] = '08:29:58'
$arr[98,1] = '09:03'
$arr[98,2] = '23_HERNANDEZ'
$arr[98,3] = 34
$arr[98,4] = 'LP1912'
$arr[99,0] = '08:52:13'
$arr[99,1] = '09:03'
$arr[99,2] = '17X38_ROMERO'
$arr[99,3] = 11
$arr[99,4] = 'LP1912'
$arr[100,0] = '08:52:13'
$arr[100,1] = '09:04'
$arr[100,2] = '23_HERNANDEZ'
$arr[100,3] = 12
$arr[100,4] = 'LP1912'
$arr[101,0] = '08:40:53'
$arr[101,1] = '09:07'
$arr[101,2] = '23_HERNANDEZ'
$arr[101,3] = 27
$arr[101,4] = 'LP1912'
$arr[102,0] = '08:10:38'
$arr[102,1] = '09:10'
$arr[102,2] = '27_EL RETIRO'
$arr[102,3] = 60
$arr[102,4] = 'LP1912'
$arr[103,0] = '07:50:33'
$arr[103,1] = '09:12'
$arr[103,2] = '27_EL RETIRO'
$arr[103,3] = 82
$arr[103,4] = 'LP1912'
$arr[104,0] = '08:40:53'
$arr[104,1] = '09:14'
$arr[104,2] = '11_ETCHEVERRY'
$arr[104,3] = 34
$arr[104,4] = 'LP1912'
$arr[105,0] = '08:52:13'
$arr[105,1] = '09:15'
$arr[105,2] = '11_ETCHEVERRY'
$arr[105,3] = 23
$arr[105,4] = 'LP1912'
$arr[106,0] = '07:38:30'
$arr[106,1] = '09:15'
$arr[106,2] = '27_EL RETIRO'
$arr[106,3] = 97
$arr[106,4] = 'LP1912'
$arr[107,0] = '08:40:53'
$arr[107,1] = '09:16'
$arr[107,2] = '27_EL RETIRO'
$arr[107,3] = 36
$arr[107,4] = 'LP1912'
$arr[108,0] = '08:52:13'
$arr[108,1] = '09:17'
$arr[108,2] = '27_EL RETIRO'
$arr[108,3] = 25
$arr[108,4] = 'LP1912'
$arr[109,0] = '08:40:53'
$arr[109,1] = '09:18'
$arr[109,2] = '215_EL PELIGRO'
$arr[109,3] = 38
$arr[109,4] = 'LP1912'
$arr[110,0] = '08:52:13'
$arr[110,1] = '09:19'
$arr[110,2] = '215_EL PELIGRO'
$arr[110,3] = 27
$arr[110,4] = 'LP1912'
$arr[111,0] = '09:22:27'
$arr[111,1] = '09:26'
$arr[111,2] = '10_OLMOS'
$arr[111,3] = 4
$arr[111,4] = 'LP1912'
$arr[112,0] = '08:29:58'
$arr[112,1] = '09:26'
$arr[112,2] = '23_HERNANDEZ'
$arr[112,3] = 57
$arr[112,4] = 'LP1912'
$arr[113,0] = '08:40:53'
$arr[113,1] = '09:28'
$arr[113,2] = '10_OLMOS'
$arr[113,3] = 48
$arr[113,4] = 'LP1912'
$arr[114,0] = '08:40:53'
$arr[114,1] = '09:29'
$arr[114,2] = '23_HERNANDEZ'
$arr[114,3] = 49
$arr[114,4] = 'LP1912'
$arr[115,0] = '08:52:13'
$arr[115,1] = '09:29'
$arr[115,2] = '10_OLMOS'
$arr[115,3] = 37
$arr[115,4] = 'LP1912'
$arr[116,0] = '09:22:27'
$arr[116,1] = '09:33'
$arr[116,2] = '23_HERNANDEZ'
$arr[116,3] = 11
$arr[116,4] = 'LP1912'
$arr[117,0] = '08:29:58'
$arr[117,1] = '09:33'
$arr[117,2] = '15_ABASTO'
$arr[117,3] = 64
$arr[117,4] = 'LP1912'
$arr[118,0] = '09:22:27'
$arr[118,1] = '09:34'
$arr[118,2] = '15_ABASTO'
$arr[118,3] = 12
$arr[118,4] = 'LP1912'
$arr[119,0] = '09:22:27'
$arr[119,1] = '09:41'
$arr[119,2] = '16_SANTA ANA'
$arr[119,3] = 19
$arr[119,4] = 'LP1912'
$arr[120,0] = '08:40:53'
$arr[120,1] = '09:44'
$arr[120,2] = '14_ABASTO'
$arr[120,3] = 64
$arr[120,4] = 'LP1912'
$arr[121,0] = '08:52:13'
$arr[121,1] = '09:45'
$arr[121,2] = '14_ABASTO'
$arr[121,3] = 53
$arr[121,4] = 'LP1912'
$arr[122,0] = '08:29:58'
$arr[122,1] = '09:48'
$arr[122,2] = '15_ABASTO'
$arr[122,3] = 79
$arr[122,4] = 'LP1912'
$arr[123,0] = '08:10:38'
$arr[123,1] = '09:49'
$arr[123,2] = '15_ABASTO'
$arr[123,3] = 99
$arr[123,4] = 'LP1912'
$arr[124,0] = '08:29:58'
$arr[124,1] = '09:50'
$arr[124,2] = '16_P MOR-SANTA ANA'
$arr[124,3] = 81
$arr[124,4] = 'LP1912'
$arr[125,0] = '09:22:27'
$arr[125,1] = '09:51'
$arr[125,2] = '16_P MOR-SANTA ANA'
$arr[125,3] = 29
$arr[125,4] = 'LP1912'
$arr[126,0] = '09:22:27'
$arr[126,1] = '09:56'
$arr[126,2] = '10_OLMOS'
$arr[126,3] = 34
$arr[126,4] = 'LP1912'
$arr[127,0] = '09:22:27'
$arr[127,1] = '10:03'
$arr[127,2] = '23_HERNANDEZ'
$arr[127,3] = 41
$arr[127,4] = 'LP1912'
$arr[128,0] = '08:40:53'
$arr[128,1] = '10:03'
$arr[128,2] = '215C_EL PATO'
$arr[128,3] = 83
$arr[128,4] = 'LP1912'
$arr[129,0] = '09:22:27'
$arr[129,1] = '10:04'
$arr[129,2] = '215C_EL PATO'
$arr[129,3] = 42
$arr[129,4] = 'LP1912'
$arr[130,0] = '09:22:27'
$arr[130,1] = '10:08'
$arr[130,2] = '11_ETCHEVERRY'
$arr[130,3] = 46
$arr[130,4] = 'LP1912'
$arr[131,0] = '10:06:07'
$arr[131,1] = '10:09'
$arr[131,2] = '11_ETCHEVERRY'
$arr[131,3] = 3
$arr[131,4] = 'LP1912'
$arr[132,0] = '10:06:07'
$arr[132,1] = '10:15'
$arr[132,2] = '16_SANTA ANA'
$arr[132,3] = 9
$arr[132,4] = 'LP1912'
$arr[133,0] = '08:40:53'
$arr[133,1] = '10:18'
$arr[133,2] = '17_ROMERO'
$arr[133,3] = 98
$arr[133,4] = 'LP1912'
$arr[134,0] = '09:22:27'
$arr[134,1] = '10:19'
$arr[134,2] = '17_ROMERO'
$arr[134,3] = 57
$arr[134,4] = 'LP1912'
$arr[135,0] = '10:06:07'
$arr[135,1] = '10:20'
$arr[135,2] = '10_OLMOS'
$arr[135,3] = 14
$arr[135,4] = 'LP1912'
$arr[136,0] = '09:22:27'
$arr[136,1] = '10:32'
$arr[136,2] = '14_ABASTO'
$arr[136,3] = 70
$arr[136,4] = 'LP1912'
$arr[137,0] = '10:06:07'
$arr[137,1] = '10:33'
$arr[137,2] = '14_ABASTO'
$arr[137,3] = 27
$arr[137,4] = 'LP1912'
$arr[138,0] = '10:06:07'
$arr[138,1] = '10:34'
$arr[138,2] = '15_ABASTO'
$arr[138,3] = 28
$arr[138,4] = 'LP1912'
$arr[139,0] = '10:06:07'
$arr[139,1] = '10:34'
$arr[139,2] = '23_HERNANDEZ'
$arr[139,3] = 28
$arr[139,4] = 'LP1912'
$arr[140,0] = '09:22:27'
$arr[140,1] = '10:40'
$arr[140,2] = '16_SANTA ANA'
$arr[140,3] = 78
$arr[140,4] = 'LP1912'
$arr[141,0] = '10:06:07'
$arr[141,1] = '10:41'
$arr[141,2] = '16_SANTA ANA'
$arr[141,3] = 35
$arr[141,4] = 'LP1912'
$arr[142,0] = '10:39:11'
$arr[142,1] = '10:44'
$arr[142,2] = '10_OLMOS'
$arr[142,3] = 5
$arr[142,4] = 'LP1912'
$arr[143,0] = '10:39:11'
$arr[143,1] = '10:49'
$arr[143,2] = '15_ABASTO'
$arr[143,3] = 10
$arr[143,4] = 'LP1912'
$arr[144,0] = '10:39:11'
$arr[144,1] = '10:51'
$arr[144,2] = '16_P MOR-SANTA ANA'
$arr[144,3] = 12
$arr[144,4] = 'LP1912'
$arr[145,0] = '10:51:29'
$arr[145,1] = '10:53'
$arr[145,2] = '16_P MOR-SANTA ANA'
$arr[145,3] = 2
$arr[145,4] = 'LP1912'
$arr[146,0] = '10:39:11'
$arr[146,1] = '10:54'
$arr[146,2] = '14_ABASTO'
$arr[146,3] = 15
$arr[146,4] = 'LP1912'
$arr[147,0] = '10:51:29'
$arr[147,1] = '10:56'
$arr[147,2] = '14_ABASTO'
$arr[147,3] = 5
$arr[147,4] = 'LP1912'
$arr[148,0] = '10:51:29'
$arr[148,1] = '10:56'
$arr[148,2] = '27_EL RETIRO'
$arr[148,3] = 5
$arr[148,4] = 'LP1912'
$arr[149,0] = '10:39:11'
$arr[149,1] = '10:57'
$arr[149,2] = '14_ABASTO'
$arr[149,3] = 18
$arr[149,4] = 'LP1912'
$arr[150,0] = '10:39:11'
$arr[150,1] = '10:57'
$arr[150,2] = '27_EL RETIRO'
$arr[150,3] = 18
$arr[150,4] = 'LP1912'
$arr[151,0] = '10:06:07'
$arr[151,1] = '10:57'
$arr[151,2] = '23_HERNANDEZ'
$arr[151,3] = 51
$arr[151,4] = 'LP1912'
$arr[152,0] = '10:39:11'
$arr[152,1] = '11:01'
$arr[152,2] = '16_SANTA ANA'
$arr[152,3] = 22
$arr[152,4] = 'LP1912'
$arr[153,0] = '10:06:07'
$arr[153,1] = '11:04'
$arr[153,2] = '17_ROMERO'
$arr[153,3] = 58
$arr[153,4] = 'LP1912'
$arr[154,0] = '10:51:29'
$arr[154,1] = '11:08'
$arr[154,2] = '225_C ROCA-H SUR'
$arr[154,3] = 17
$arr[154,4] = 'LP1912'
$arr[155,0] = '10:51:29'
$arr[155,1] = '11:09'
$arr[155,2] = '17_ROMERO'
$arr[155,3] = 18
$arr[155,4] = 'LP1912'
$arr[156,0] = '09:22:27'
$arr[156,1] = '11:09'
$arr[156,2] = '14_ABASTO'
$arr[156,3] = 107
$arr[156,4] = 'LP1912'
$arr[157,0] = '10:39:11'
$arr[157,1] = '11:11'
$arr[157,2] = '17_ROMERO'
$arr[157,3] = 32
$arr[157,4] = 'LP1912'
$arr[158,0] = '10:51:29'
$arr[158,1] = '11:19'
$arr[158,2] = '215C_EL PATO'
$arr[158,3] = 28
$arr[158,4] = 'LP1912'
$arr[159,0] = '10:51:29'
$arr[159,1] = '11:20'
$arr[159,2] = '11_ETCHEVERRY'
$arr[159,3] = 29
$arr[159,4] = 'LP1912'
$arr[160,0] = '10:39:11'
$arr[160,1] = '11:21'
$arr[160,2] = '11_ETCHEVERRY'
$arr[160,3] = 42
$arr[160,4] = 'LP1912'
$arr[161,0] = '10:39:11'
$arr[161,1] = '11:21'
$arr[161,2] = '23_HERNANDEZ'
$arr[161,3] = 42
$arr[161,4] = 'LP1912'
$arr[162,0] = '10:39:11'
$arr[162,1] = '11:30'
$arr[162,2] = '23_HERNANDEZ'
$arr[162,3] = 51
$arr[162,4] = 'LP1912'
$arr[163,0] = '10:06:07'
$arr[163,1] = '11:30'
$arr[163,2] = '14_ABASTO'
$arr[163,3] = 84
$arr[163,4] = 'LP1912'
$arr[164,0] = '10:51:29'
$arr[164,1] = '11:32'
$arr[164,2] = '10_OLMOS'
$arr[164,3] = 41
$arr[164,4] = 'LP1912'
$arr[165,0] = '10:51:29'
$arr[165,1] = '11:33'
$arr[165,2] = '215A_EL PATO'
$arr[165,3] = 42
$arr[165,4] = 'LP1912'
$arr[166,0] = '10:51:29'
$arr[166,1] = '11:33'
$arr[166,2] = '23_HERNANDEZ'
$arr[166,3] = 42
$arr[166,4] = 'LP1912'
$arr[167,0] = '10:51:29'
$arr[167,1] = '11:41'
$arr[167,2] = '16_SANTA ANA'
$arr[167,3] = 50
$arr[167,4] = 'LP1912'
$arr[168,0] = '10:51:29'
$arr[168,1] = '11:44'
$arr[168,2] = '215B_EL PATO'
$arr[168,3] = 53
$arr[168,4] = 'LP1912'
$arr[169,0] = '10:39:11'
$arr[169,1] = '11:45'
$arr[169,2] = '215B_EL PATO'
$arr[169,3] = 66
$arr[169,4] = 'LP1912'
$arr[170,0] = '10:51:29'
$arr[170,1] = '11:49'
$arr[170,2] = '15_ABASTO'
$arr[170,3] = 58
$arr[170,4] = 'LP1912'
$arr[171,0] = '10:51:29'
$arr[171,1] = '11:51'
$arr[171,2] = '16_P MOR-SANTA ANA'
$arr[171,3] = 60
$arr[171,4] = 'LP1912'
$arr[172,0] = '10:51:29'
$arr[172,1] = '11:56'
$arr[172,2] = '225_GOMEZ'
$arr[172,3] = 65
$arr[172,4] = 'LP1912'
$arr[173,0] = '10:51:29'
$arr[173,1] = '11:56'
$arr[173,2] = '16_SANTA ANA'
$arr[173,3] = 65
$arr[173,4] = 'LP1912'
$arr[174,0] = '10:51:29'
$arr[174,1] = '12:04'
$arr[174,2] = '17_ROMERO'
$arr[174,3] = 73
$arr[174,4] = 'LP1912'
$arr[175,0] = '10:51:29'
$arr[175,1] = '12:08'
$arr[175,2] = '14_ABASTO'
$arr[175,3] = 77
$arr[175,4] = 'LP1912'
$arr[176,0] = '10:51:29'
$arr[176,1] = '12:19'
$arr[176,2] = '15_ABASTO'
$arr[176,3] = 88
$arr[176,4] = 'LP1912'
$arr[177,0] = '10:51:29'
$arr[177,1] = '12:20'
$arr[177,2] = '10_OLMOS'
$arr[177,3] = 89
$arr[177,4] = 'LP1912'
$arr[178,0] = '10:51:29'
$arr[178,1] = '12:32'
$arr[178,2] = '11_ETCHEVERRY'
$arr[178,3] = 101
$arr[178,4] = 'LP1912'
$arr[179,0] = '10:39:11'
$arr[179,1] = '12:33'
$arr[179,2] = '11_ETCHEVERRY'
$arr[179,3] = 114
$arr[179,4] = 'LP1912'
$arr[180,0] = '10:51:29'
$arr[180,1] = '12:34'
$arr[180,2] = '215C_EL PATO'
$arr[180,3] = 103
$arr[180,4] = 'LP1912'
$arr[181,0] = '10:51:29'
$arr[181,1] = '12:36'
$arr[181,2] = '27_EL RETIRO'
$arr[181,3] = 105
$arr[181,4] = 'LP1912'
$arr[182,0] = '10:39:11'
$arr[182,1] = '12:37'
$arr[182,2] = '27_EL RETIRO'
$arr[182,3] = 118
$arr[182,4] = 'LP1912'
$ws.Range("A5:E187").Value = $arr

# ===== Sheet 2 =====
$ws = $wb.Worksheets.Item(2)

$ws.Cells.Item(1, 1).Value = 'LÍNEA 141 - LP1912-215 - 24/01/2026'
$ws.Cells.Item(2, 1).Value = 'Última actualización: 10:51:29'
$ws.Cells.Item(3, 1).Value = 'Total filas: 27'

$arr = New-Object "object[,]" 28,5
$arr[0,0] = 'Hora_Scrap'
$arr[0,1] = 'Hora_Llegada'
$arr[0,2] = 'Linea'
$arr[0,3] = 'Minutos'
$arr[0,4] = 'Parada'
$arr[1,0] = '04:44:46'
$arr[1,1] = '04:46'
$arr[1,2] = '215_EL PELIGRO'
$arr[1,3] = 2
$arr[1,4] = 'LP1912'
$arr[2,0] = '03:52:04'
$arr[2,1] = '04:46'
$arr[2,2] = '215A_EL PATO'
$arr[2,3] = 54
$arr[2,4] = 'LP1912'
$arr[3,0] = '04:32:18'
$arr[3,1] = '04:47'
$arr[3,2] = '215_EL PELIGRO'
$arr[3,3] = 15
$arr[3,4] = 'LP1912'
$arr[4,0] = '03:52:04'
$arr[4,1] = '05:35'
$arr[4,2] = '215B_EL PATO'
$arr[4,3] = 103
$arr[4,4] = 'LP1912'
$arr[5,0] = '05:57:38'
$arr[5,1] = '06:10'
$arr[5,2] = '215A_EL PATO'
$arr[5,3] = 13
$arr[5,4] = 'LP1912'
$arr[6,0] = '05:43:29'
$arr[6,1] = '06:11'
$arr[6,2] = '215A_EL PATO'
$arr[6,3] = 28
$arr[6,4] = 'LP1912'
$arr[7,0] = '06:46:37'
$arr[7,1] = '06:50'
$arr[7,2] = '215A_EL PATO'
$arr[7,3] = 4
$arr[7,4] = 'LP1912'
$arr[8,0] = '06:35:33'
$arr[8,1] = '06:51'
$arr[8,2] = '215A_EL PATO'
$arr[8,3] = 16
$arr[8,4] = 'LP1912'
$arr[9,0] = '06:53:56'
$arr[9,1] = '07:06'
$arr[9,2] = '215C_EL PATO'
$arr[9,3] = 13
$arr[9,4] = 'LP1912'
$arr[10,0] = '06:18:01'
$arr[10,1] = '07:07'
$arr[10,2] = '215C_EL PATO'
$arr[10,3] = 49
$arr[10,4] = 'LP1912'
$arr[11,0] = '07:12:47'
$arr[11,1] = '07:20'
$arr[11,2] = '215A_EL PATO'
$arr[11,3] = 8
$arr[11,4] = 'LP1912'
$arr[12,0] = '06:35:33'
$arr[12,1] = '07:21'
$arr[12,2] = '215A_EL PATO'
$arr[12,3] = 46
$arr[12,4] = 'LP1912'
$arr[13,0] = '08:29:58'
$arr[13,1] = '08:33'
$arr[13,2] = '215C_EL PATO'
$arr[13,3] = 4
$arr[13,4] = 'LP1912'
$arr[14,0] = '07:50:33'
$arr[14,1] = '08:34'
$arr[14,2] = '215C_EL PATO'
$arr[14,3] = 44
$arr[14,4] = 'LP1912'
$arr[15,0] = '08:40:53'
$arr[15,1] = '08:47'
$arr[15,2] = '215A_EL PATO'
$arr[15,3] = 7
$arr[15,4] = 'LP1912'
$arr[16,0] = '08:10:38'
$arr[16,1] = '08:48'
$arr[16,2] = '215A_EL PATO'
$arr[16,3] = 38
$arr[16,4] = 'LP1912'
$arr[17,0] = '08:40:53'
$arr[17,1] = '08:59'
$arr[17,2] = '215B_EL PATO'
$arr[17,3] = 19
$arr[17,4] = 'LP1912'
$arr[18,0] = '08:52:13'
$arr[18,1] = '09:00'
$arr[18,2] = '215B_EL PATO'
$arr[18,3] = 8
$arr[18,4] = 'LP1912'
$arr[19,0] = '08:40:53'
$arr[19,1] = '09:18'
$arr[19,2] = '215_EL PELIGRO'
$arr[19,3] = 38
$arr[19,4] = 'LP1912'
$arr[20,0] = '08:52:13'
$arr[20,1] = '09:19'
$arr[20,2] = '215_EL PELIGRO'
$arr[20,3] = 27
$arr[20,4] = 'LP1912'
$arr[21,0] = '08:40:53'
$arr[21,1] = '10:03'
$arr[21,2] = '215C_EL PATO'
$arr[21,3] = 83
$arr[21,4] = 'LP1912'
$arr[22,0] = '09:22:27'
$arr[22,1] = '10:04'
$arr[22,2] = '215C_EL PATO'
$arr[22,3] = 42
$arr[22,4] = 'LP1912'
$arr[23,0] = '10:51:29'
$arr[23,1] = '11:19'
$arr[23,2] = '215C_EL PATO'
$arr[23,3] = 28
$arr[23,4] = 'LP1912'
$arr[24,0] = '10:51:29'
$arr[24,1] = '11:33'
$arr[24,2] = '215A_EL PATO'
$arr[24,3] = 42
$arr[24,4] = 'LP1912'
$arr[25,0] = '10:51:29'
$arr[25,1] = '11:44'
$arr[25,2] = '215B_EL PATO'
$arr[25,3] = 53
$arr[25,4] = 'LP1912'
$arr[26,0] = '10:39:11'
$arr[26,1] = '11:45'
$arr[26,2] = '215B_EL PATO'
$arr[26,3] = 66
$arr[26,4] = 'LP1912'
$arr[27,0] = '10:51:29'
$arr[27,1] = '12:34'
$arr[27,2] = '215C_EL PATO'
$arr[27,3] = 103
$arr[27,4] = 'LP1912'
$ws.Range("A5:E32").Value = $arr

# ===== Sheet 3 =====
$ws = $wb.Worksheets.Item(3)

$ws.Cells.Item(1, 1).Value = 'LÍNEA 141 - 6203-6173 - 24/01/2026'
$ws.Cells.Item(2, 1).Value = 'Última actualización: 10:51:29'
$ws.Cells.Item(3, 1).Value = 'Total filas: 17'

$arr = New-Object "object[,]" 18,5
$arr[0,0] = 'Hora_Scrap'
$arr[0,1] = 'Hora_Llegada'
$arr[0,2] = 'Linea'
$arr[0,3] = 'Minutos'
$arr[0,4] = 'Parada'
$arr[1,0] = '03:52:04'
$arr[1,1] = '05:44'
$arr[1,2] = '215A_LA PLATA'
$arr[1,3] = 112
$arr[1,4] = 'L6173'
$arr[2,0] = '07:12:47'
$arr[2,1] = '07:27'
$arr[2,2] = '215A_LA PLATA'
$arr[2,3] = 15
$arr[2,4] = 'L6173'
$arr[3,0] = '07:12:47'
$arr[3,1] = '08:09'
$arr[3,2] = '215A_LA PLATA'
$arr[3,3] = 57
$arr[3,4] = 'L6173'
$arr[4,0] = '07:50:33'
$arr[4,1] = '08:10'
$arr[4,2] = '215A_LA PLATA'
$arr[4,3] = 20
$arr[4,4] = 'L6173'
$arr[5,0] = '08:10:38'
$arr[5,1] = '08:12'
$arr[5,2] = '215A_LA PLATA'
$arr[5,3] = 2
$arr[5,4] = 'L6173'
$arr[6,0] = '07:12:47'
$arr[6,1] = '08:22'
$arr[6,2] = '215C_LA PLATA'
$arr[6,3] = 70
$arr[6,4] = 'L6203'
$arr[7,0] = '07:50:33'
$arr[7,1] = '08:23'
$arr[7,2] = '215C_LA PLATA'
$arr[7,3] = 33
$arr[7,4] = 'L6203'
$arr[8,0] = '08:10:38'
$arr[8,1] = '08:24'
$arr[8,2] = '215C_LA PLATA'
$arr[8,3] = 14
$arr[8,4] = 'L6203'
$arr[9,0] = '08:40:53'
$arr[9,1] = '08:51'
$arr[9,2] = '215A_LA PLATA'
$arr[9,3] = 11
$arr[9,4] = 'L6173'
$arr[10,0] = '08:52:13'
$arr[10,1] = '08:52'
$arr[10,2] = '215A_LA PLATA'
$arr[10,3] = 0
$arr[10,4] = 'L6173'
$arr[11,0] = '08:40:53'
$arr[11,1] = '09:55'
$arr[11,2] = '215C_LA PLATA'
$arr[11,3] = 75
$arr[11,4] = 'L6203'
$arr[12,0] = '09:22:27'
$arr[12,1] = '09:56'
$arr[12,2] = '215C_LA PLATA'
$arr[12,3] = 34
$arr[12,4] = 'L6203'
$arr[13,0] = '08:40:53'
$arr[13,1] = '10:10'
$arr[13,2] = '215A_LA PLATA'
$arr[13,3] = 90
$arr[13,4] = 'L6173'
$arr[14,0] = '10:06:07'
$arr[14,1] = '10:11'
$arr[14,2] = '215A_LA PLATA'
$arr[14,3] = 5
$arr[14,4] = 'L6173'
$arr[15,0] = '08:40:53'
$arr[15,1] = '10:21'
$arr[15,2] = '215B_LP-P MOR-1 Y 57'
$arr[15,3] = 101
$arr[15,4] = 'L6173'
$arr[16,0] = '10:06:07'
$arr[16,1] = '10:22'
$arr[16,2] = '215B_LP-P MOR-1 Y 57'
$arr[16,3] = 16
$arr[16,4] = 'L6173'
$arr[17,0] = '10:51:29'
$arr[17,1] = '11:56'
$arr[17,2] = '215C_LA PLATA'
$arr[17,3] = 65
$arr[17,4] = 'L6203'
$ws.Range("A5:E22").Value = $arr
